# Update the "Time:" timestamp embedded in the statsmodels OLS summary text
# block that lives in cell B2 of every worksheet (one per backward-
# elimination step, tabs "41" down to "14"). Only the wall-clock time the
# summary was generated changes (20:51:50 -> 20:59:52); every other
# character of the text, including the unchanged "Date:" line, is left
# untouched.

$wb = $excel.ActiveWorkbook

$oldTime = "20:51:50"
$newTime = "20:59:52"

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("B2")
    $text = $cell.Text
    if ($text -and $text.Contains($oldTime)) {
        $cell.Value = $text.Replace($oldTime, $newTime)
    }
}
